$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly completed task-days for week 11 (col K) and week 12 (col L)
$ws.Range("K24").Value = 9
$ws.Range("K27").Value = 9
$ws.Range("K28").Value = 9
$ws.Range("L28").Value = 9
$ws.Range("K29").Value = 9
$ws.Range("L30").Value = 9

# Extend the "% complete" formula in column M from C:J to C:L so it
# picks up the newly-filled week 11 / week 12 columns, and refill the
# shared formula down through row 30.
$ws.Range("M17:M30").Formula = "=SUM(C17:L17)/B17*100"

# Match the selection left by the author's editing session
$ws.Range("C17:L30").Select()

$wb.Save()
